# "Update to Mod 3"
# Title slide (slide 1) subtitle text is refreshed to describe the
# course's new Linux + Docker focus (previously referenced Azure App
# Service's Linux offerings specifically).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$subtitle = $null
foreach ($shp in $s.Shapes) {
    if ($shp.Id -eq 3) {
        $subtitle = $shp
    }
}
if ($subtitle -eq $null) {
    $subtitle = $s.Shapes.Item(2)
}

$subtitle.TextFrame.TextRange.Text = "Learning Linux and Docker for Anyone"
